$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh (GitHub Actions): update price (D) and
# 1h-volume-change (E) columns for each ranked coin; rows 49/50 also
# swapped rank order (SuiNetwork now above ONDO).

$ws.Range("D2").Value = "69.422.98"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "3.692.69"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'692.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'162.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.37%  "
$ws.Range("D7").Value = "3.689.59"
$ws.Range("E7").Value = "  -3.07%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -4.68%  "
$ws.Range("D10").Value = "'0.147"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -8.38%  "
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "'0.442"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.33%  "
$ws.Range("D13").Value = "'0.0000238"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.25%  "
$ws.Range("D14").Value = "'33.31"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.17%  "
$ws.Range("D15").Value = "4.314.75"
$ws.Range("E15").Value = "  -3.05%  "
$ws.Range("D16").Value = "3.690.18"
$ws.Range("E16").Value = "  -3.59%  "
$ws.Range("D17").Value = "69.477.79"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").Value = "'16.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -7.76%  "
$ws.Range("D20").Value = "'6.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -8.03%  "
$ws.Range("D21").Value = "'479.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.78%  "
$ws.Range("D22").Value = "'9.99"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.75%  "
$ws.Range("D23").Value = "'0.662"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.30%  "
$ws.Range("D24").Value = "'79.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.77%  "
$ws.Range("D25").Value = "3.838.27"
$ws.Range("D26").Value = "'0.0000129"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -9.37%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'11.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.03%  "
$ws.Range("D29").Value = "'9.46"
$ws.Range("D29").ClearFormats()
$ws.Range("E30").Value = "  -11.24%  "
$ws.Range("E31").Value = "  -9.97%  "
$ws.Range("D32").Value = "'6.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.62%  "
$ws.Range("E33").Value = "  -7.54%  "
$ws.Range("E34").Value = "  -5.26%  "
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "'26.93"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.20%  "
$ws.Range("D37").Value = "3.657.30"
$ws.Range("E37").Value = "  -3.00%  "
$ws.Range("D38").Value = "'8.44"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.40%  "
$ws.Range("D39").Value = "'6.26"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("E41").Value = "  -8.21%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -6.37%  "
$ws.Range("D45").Value = "'163.52"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.07%  "
$ws.Range("D46").Value = "'48.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").Value = "'30.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  -14.93%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").Value = "'1.14"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.34"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("E51").Value = "  -9.00%  "
